$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two period-end dates (shift each by 10 days: 2016-12-25 -> 2017-01-04, 2016-12-29 -> 2017-01-08)
$ws.Range("B1").Value = 42739
$ws.Range("D1").Value = 42743

# Clear out the now-unused "weights" helper column (B) for the data rows.
$ws.Range("B3:B30").ClearContents()
